$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 2923.75
$ws.Range("I16").Value = 2923.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2923.75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2693.75
$ws.Range("N16").ClearContents()

$ws.Range("H19").Value = 3543.7812
$ws.Range("I19").Value = 7043.933
$ws.Range("J19").Value = 455.41177
$ws.Range("K19").Value = 7043.933
$ws.Range("L19").Value = 455.41177
$ws.Range("M19").Value = -6868.933
$ws.Range("N19").Value = -805.4117699999999

$ws.Range("H32").Value = 7744573
$ws.Range("I32").Value = 363
$ws.Range("J32").Value = 11616678
$ws.Range("K32").Value = 363
$ws.Range("L32").Value = 11616678
$ws.Range("M32").Value = -37
$ws.Range("N32").Value = -11617330

$ws.Range("H113").Value = 3631.4827
$ws.Range("I113").Value = 3196.6365
$ws.Range("J113").Value = 3897.2222
$ws.Range("K113").Value = 3196.6365
$ws.Range("L113").Value = 3897.2222
$ws.Range("M113").Value = 57.36349999999993
$ws.Range("N113").Value = -10405.2222

$ws.Range("H116").Value = 114230.79
$ws.Range("I116").Value = 194062.27
$ws.Range("J116").Value = 4462.5
$ws.Range("K116").Value = 194062.27
$ws.Range("L116").Value = 4462.5
$ws.Range("M116").Value = -190620.27
$ws.Range("N116").Value = -11346.5

$ws.Range("H125").Value = 941.6667
$ws.Range("I125").Value = 1086.7778
$ws.Range("J125").Value = 724
$ws.Range("K125").Value = 9781.0002
$ws.Range("L125").Value = 6516
$ws.Range("M125").Value = -7321.0002
$ws.Range("N125").Value = -11436

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 903
$ws.Range("I2").Value = 874.2
$ws.Range("K2").Value = 874.2
$ws.Range("M2").Value = -761.2

$ws.Range("H33").Value = 5756.5
$ws.Range("I33").Value = 26
$ws.Range("J33").Value = 7666.6665
$ws.Range("K33").Value = 26
$ws.Range("L33").Value = 7666.6665
$ws.Range("M33").Value = 303
$ws.Range("N33").Value = -8324.666499999999

$ws.Range("H45").Value = 1032.8276
$ws.Range("I45").Value = 870.1053000000001
$ws.Range("J45").Value = 1342
$ws.Range("K45").Value = 870.1053000000001
$ws.Range("L45").Value = 1342
$ws.Range("M45").Value = -493.1053000000001
$ws.Range("N45").Value = -2096

$ws.Range("H63").Value = 4764.6924
$ws.Range("I63").Value = 5284.15
$ws.Range("J63").Value = 3033.1667
$ws.Range("K63").Value = 5284.15
$ws.Range("L63").Value = 3033.1667
$ws.Range("M63").Value = -4598.15
$ws.Range("N63").Value = -4405.1667

$ws.Range("H66").Value = 4764.6924
$ws.Range("I66").Value = 5284.15
$ws.Range("J66").Value = 3033.1667
$ws.Range("K66").Value = 26420.75
$ws.Range("L66").Value = 15165.8335
$ws.Range("M66").Value = -22988.75
$ws.Range("N66").Value = -22029.8335

$ws.Range("H74").Value = 1781.8572
$ws.Range("I74").Value = 1062.55
$ws.Range("J74").Value = 3580.125
$ws.Range("K74").Value = 1062.55
$ws.Range("L74").Value = 3580.125
$ws.Range("M74").Value = -188.55
$ws.Range("N74").Value = -5328.125

$ws.Range("H77").Value = 1781.8572
$ws.Range("I77").Value = 1062.55
$ws.Range("J77").Value = 3580.125
$ws.Range("K77").Value = 5312.75
$ws.Range("L77").Value = 17900.625
$ws.Range("M77").Value = -944.75
$ws.Range("N77").Value = -26636.625

$ws.Range("H110").Value = 1401.8125
$ws.Range("I110").Value = 1562.1666
$ws.Range("J110").Value = 1195.6428
$ws.Range("K110").Value = 1562.1666
$ws.Range("L110").Value = 1195.6428
$ws.Range("M110").Value = 482.8334
$ws.Range("N110").Value = -5285.6428

$ws.Range("H116").Value = 903
$ws.Range("I116").Value = 874.2
$ws.Range("K116").Value = 874.2
$ws.Range("M116").Value = 1419.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 903
$ws.Range("I3").Value = 874.2
$ws.Range("K3").Value = 874.2
$ws.Range("M3").Value = -760.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2831.5774
$ws.Range("I31").Value = 1858.1951
$ws.Range("J31").Value = 4161.8667
$ws.Range("K31").Value = 1858.1951
$ws.Range("L31").Value = 4161.8667
$ws.Range("M31").Value = -1563.1951
$ws.Range("N31").Value = -4751.8667

$ws.Range("H34").Value = 2831.5774
$ws.Range("I34").Value = 1858.1951
$ws.Range("J34").Value = 4161.8667
$ws.Range("K34").Value = 1858.1951
$ws.Range("L34").Value = 4161.8667
$ws.Range("M34").Value = -1656.1951
$ws.Range("N34").Value = -4565.8667

$ws.Range("H62").Value = 3169.3696
$ws.Range("I62").Value = 2967.818
$ws.Range("K62").Value = 2967.818
$ws.Range("M62").Value = -2343.818

$ws.Range("H65").Value = 3169.3696
$ws.Range("I65").Value = 2967.818
$ws.Range("K65").Value = 14839.09
$ws.Range("M65").Value = -11719.09

$ws.Range("H94").Value = 7425.077
$ws.Range("I94").Value = 1242.4
$ws.Range("J94").Value = 11289.25
$ws.Range("K94").Value = 1242.4
$ws.Range("L94").Value = 11289.25
$ws.Range("M94").Value = -791.4000000000001
$ws.Range("N94").Value = -12191.25

$ws.Range("H99").Value = 55318.156
$ws.Range("I99").Value = 64648.688
$ws.Range("J99").Value = 5555.3335
$ws.Range("K99").Value = 64648.688
$ws.Range("L99").Value = 5555.3335
$ws.Range("M99").Value = -63150.688
$ws.Range("N99").Value = -8551.333500000001

$ws.Range("H126").Value = 55318.156
$ws.Range("I126").Value = 64648.688
$ws.Range("J126").Value = 5555.3335
$ws.Range("K126").Value = 193946.064
$ws.Range("L126").Value = 16666.0005
$ws.Range("M126").Value = -191476.064
$ws.Range("N126").Value = -21606.0005

$ws.Range("H132").Value = 2187.6584
$ws.Range("I132").Value = 1052.4546
$ws.Range("J132").Value = 3502.1052
$ws.Range("K132").Value = 3157.3638
$ws.Range("L132").Value = 10506.3156
$ws.Range("M132").Value = -627.3638000000001
$ws.Range("N132").Value = -15566.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 36.8
$ws.Range("I12").Value = 23.666666
$ws.Range("K12").Value = 70.99999800000001
$ws.Range("M12").Value = 102.000002

$ws.Range("H14").Value = 652.56525
$ws.Range("I14").Value = 652.56525
$ws.Range("K14").Value = 1957.69575
$ws.Range("M14").Value = -1784.69575

$ws.Range("H34").Value = 904.04
$ws.Range("J34").Value = 952.2174
$ws.Range("L34").Value = 2856.6522
$ws.Range("N34").Value = -3024.6522

$ws.Range("H68").Value = 571.7143
$ws.Range("I68").Value = 502
$ws.Range("J68").Value = 583.3333
$ws.Range("K68").Value = 1506
$ws.Range("L68").Value = 1749.9999
$ws.Range("M68").Value = -695
$ws.Range("N68").Value = -3371.9999

$ws.Range("H71").Value = 571.7143
$ws.Range("I71").Value = 502
$ws.Range("J71").Value = 583.3333
$ws.Range("K71").Value = 4518
$ws.Range("L71").Value = 5249.9997
$ws.Range("M71").Value = -462
$ws.Range("N71").Value = -13361.9997

$ws.Range("H80").Value = 1422
$ws.Range("J80").Value = 1487.25
$ws.Range("L80").Value = 4461.75
$ws.Range("N80").Value = -6333.75

$ws.Range("H83").Value = 1422
$ws.Range("J83").Value = 1487.25
$ws.Range("L83").Value = 13385.25
$ws.Range("N83").Value = -22745.25

$ws.Range("H92").Value = 901.5
$ws.Range("J92").Value = 902.625
$ws.Range("L92").Value = 2707.875
$ws.Range("N92").Value = -5203.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 4000
$ws.Range("I27").Value = 4000
$ws.Range("K27").Value = 4000
$ws.Range("M27").Value = -3834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 625
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 625
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 625
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1215

$ws.Range("H27").Value = 625
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 625
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 625
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -839

$ws.Range("H55").Value = 211.03847
$ws.Range("I55").Value = 171.1579
$ws.Range("J55").Value = 319.2857
$ws.Range("K55").Value = 171.1579
$ws.Range("L55").Value = 319.2857
$ws.Range("M55").Value = 1.842099999999988
$ws.Range("N55").Value = -665.2857

$ws.Range("H93").Value = 1477.2142
$ws.Range("I93").Value = 1566.091
$ws.Range("J93").Value = 1151.3334
$ws.Range("K93").Value = 1566.091
$ws.Range("L93").Value = 1151.3334
$ws.Range("M93").Value = -318.0909999999999
$ws.Range("N93").Value = -3647.3334

$ws.Range("H122").Value = 3066.4
$ws.Range("I122").Value = 3166.2222
$ws.Range("J122").Value = 2916.6667
$ws.Range("K122").Value = 9498.6666
$ws.Range("L122").Value = 8750.000100000001
$ws.Range("M122").Value = -7048.6666
$ws.Range("N122").Value = -13650.0001
